$d = $word.ActiveDocument

# Ordered list of (old, new) replacement pairs, exactly matching the
# left-to-right, top-to-bottom order of the division problems in the table.
$replacements = @(
    "98÷6=", "32÷9=",
    "34÷9=", "81÷2=",
    "72÷8=", "23÷7=",
    "10÷4=", "65÷8=",
    "73÷5=", "23÷4=",
    "64÷3=", "54÷3=",
    "48÷9=", "88÷9=",
    "93÷6=", "45÷2=",
    "44÷5=", "34÷9=",
    "81÷2=", "66÷2=",
    "16÷3=", "82÷9=",
    "78÷8=", "21÷6=",
    "98÷8=", "33÷9=",
    "33÷6=", "90÷7=",
    "33÷8=", "75÷5=",
    "75÷3=", "96÷9=",
    "92÷9=", "34÷4=",
    "23÷4=", "32÷4=",
    "97÷3=", "36÷5=",
    "24÷9=", "82÷8=",
    "17÷3=", "68÷8=",
    "92÷8=", "73÷3=",
    "87÷2=", "88÷8=",
    "53÷3=", "34÷2=",
    "10÷3=", "32÷8="
)

$table = $d.Tables.Item(1)

$idx = 0
for ($r = 1; $r -le $table.Rows.Count; $r++) {
    for ($c = 1; $c -le $table.Columns.Count; $c++) {
        $cell = $table.Cell($r, $c)
        $cellRange = $cell.Range
        $cellText = $cellRange.Text
        # Cell text includes trailing cell-mark characters; trim them off.
        $trimmed = $cellText.TrimEnd([char]7, [char]13)
        if ($trimmed.Length -gt 0) {
            $old = $replacements[$idx]
            $new = $replacements[$idx + 1]
            $idx += 2
            if ($trimmed -ne $old) {
                throw "Mismatch at row $r col $c : expected '$old' but found '$trimmed'"
            }
            $findRange = $cell.Range
            $findRange.End = $findRange.End - 1
            $findRange.Text = $new
        }
    }
}
